# Apply updates to column F (dSF) for several rows, per commit:
# "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -7
$ws.Range("F3").Value = -6
$ws.Range("F5").Value = -8
$ws.Range("F7").Value = -8
$ws.Range("F10").Value = 0
$ws.Range("F12").Value = -2
$ws.Range("F13").Value = 1
$ws.Range("F14").Value = -1
$ws.Range("F21").Value = 7
